$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
$ws.Columns.Item(3).ColumnWidth = 52.16666666666667
$ws.Columns.Item(4).ColumnWidth = 50.16666666666667
$ws.Columns.Item(8).ColumnWidth = 35.16666666666667

# --- Row data updates (rows 2-12 changed, rows 13-18 new) ---
# Row 2
$ws.Cells.Item(2,1).Value = '1331082'
$ws.Cells.Item(2,2).Value = 'https://aiesec.org/opportunity/global-talent/1331082'
$ws.Cells.Item(2,3).Value = 'Taste Hungary|Junior Field Service Associate'
$ws.Cells.Item(2,4).Value = 'Budapeste, Hungria'
$ws.Cells.Item(2,5).Value = 'Yes'
$ws.Cells.Item(2,6).Value = '2 applicants'
$ws.Cells.Item(2,7).Value = '6 - 18 Months'
$ws.Cells.Item(2,8).Value = 'EATON'
$ws.Cells.Item(2,5).Interior.Color = 65535

# Row 3
$ws.Cells.Item(3,1).Value = '1331094'
$ws.Cells.Item(3,2).Value = 'https://aiesec.org/opportunity/global-talent/1331094'
$ws.Cells.Item(3,3).Value = '[GBS] Transportation Sourcing Assistant'
$ws.Cells.Item(3,4).Value = 'Charles-de-Gaulle-Straße 20, 53113 Bonn, Germany'
$ws.Cells.Item(3,5).Value = 'Yes'
$ws.Cells.Item(3,6).Value = '4 applicants'
$ws.Cells.Item(3,7).Value = '6 - 18 Months'
$ws.Cells.Item(3,8).Value = 'DHL Group'
$ws.Cells.Item(3,5).Interior.Color = 65535

# Row 4
$ws.Cells.Item(4,1).Value = '1331078'
$ws.Cells.Item(4,2).Value = 'https://aiesec.org/opportunity/global-talent/1331078'
$ws.Cells.Item(4,3).Value = 'Videography Intern – Cinematic Social Media'
$ws.Cells.Item(4,4).Value = 'Karachi, Pakistan'
$ws.Cells.Item(4,5).Value = 'No'
$ws.Cells.Item(4,6).Value = '0 applicants'
$ws.Cells.Item(4,7).Value = '3 - 6 Months'
$ws.Cells.Item(4,8).Value = 'Reload snacks pvt ltd'

# Row 5
$ws.Cells.Item(5,1).Value = '1330641'
$ws.Cells.Item(5,2).Value = 'https://aiesec.org/opportunity/global-talent/1330641'
$ws.Cells.Item(5,3).Value = 'Travel Guide/Consultant'
$ws.Cells.Item(5,4).Value = 'Mehrauli, Delhi, India'
$ws.Cells.Item(5,5).Value = 'No'
$ws.Cells.Item(5,6).Value = '0 applicants'
$ws.Cells.Item(5,7).Value = '6 - 18 Months'
$ws.Cells.Item(5,8).Value = 'GeTS Holidays Private Limited'

# Row 6
$ws.Cells.Item(6,1).Value = '1328685'
$ws.Cells.Item(6,2).Value = 'https://aiesec.org/opportunity/global-talent/1328685'
$ws.Cells.Item(6,3).Value = 'Medical Advisor (Russian Speaker)'
$ws.Cells.Item(6,4).Value = 'İstanbul, Türkiye'
$ws.Cells.Item(6,5).Value = 'No'
$ws.Cells.Item(6,6).Value = '5 applicants'
$ws.Cells.Item(6,7).Value = '6 - 18 Months'
$ws.Cells.Item(6,8).Value = 'International Plus'

# Row 7
$ws.Cells.Item(7,1).Value = '1325297'
$ws.Cells.Item(7,2).Value = 'https://aiesec.org/opportunity/global-talent/1325297'
$ws.Cells.Item(7,3).Value = 'International Sales Representetive Spanish Speaker'
$ws.Cells.Item(7,4).Value = 'Maslak, Sarıyer/İstanbul, Türkiye'
$ws.Cells.Item(7,5).Value = 'No'
$ws.Cells.Item(7,6).Value = '24 applicants'
$ws.Cells.Item(7,7).Value = '6 - 18 Months'
$ws.Cells.Item(7,8).Value = 'Esvita Clinic'

# Row 8
$ws.Cells.Item(8,1).Value = '1321910'
$ws.Cells.Item(8,2).Value = 'https://aiesec.org/opportunity/global-talent/1321910'
$ws.Cells.Item(8,3).Value = 'Marketing'
$ws.Cells.Item(8,4).Value = 'Gaziantep, Türkiye'
$ws.Cells.Item(8,5).Value = 'No'
$ws.Cells.Item(8,6).Value = '117 applicants'
$ws.Cells.Item(8,7).Value = '6 - 18 Months'
$ws.Cells.Item(8,8).Value = 'Eman Agro Gıda'

# Row 9
$ws.Cells.Item(9,1).Value = '1321055'
$ws.Cells.Item(9,2).Value = 'https://aiesec.org/opportunity/global-talent/1321055'
$ws.Cells.Item(9,3).Value = 'International Sales Representetive Italian Speaker'
$ws.Cells.Item(9,4).Value = 'Maslak, Sarıyer/İstanbul, Türkiye'
$ws.Cells.Item(9,5).Value = 'No'
$ws.Cells.Item(9,6).Value = '15 applicants'
$ws.Cells.Item(9,7).Value = '6 - 18 Months'
$ws.Cells.Item(9,8).Value = 'Esvita Clinic'

# Row 10
$ws.Cells.Item(10,1).Value = '1321054'
$ws.Cells.Item(10,2).Value = 'https://aiesec.org/opportunity/global-talent/1321054'
$ws.Cells.Item(10,3).Value = 'International Sales Representetive Russian Speaker'
$ws.Cells.Item(10,4).Value = 'Maslak, Sarıyer/İstanbul, Türkiye'
$ws.Cells.Item(10,5).Value = 'No'
$ws.Cells.Item(10,6).Value = '15 applicants'
$ws.Cells.Item(10,7).Value = '6 - 18 Months'
$ws.Cells.Item(10,8).Value = 'Esvita Clinic'

# Row 11
$ws.Cells.Item(11,1).Value = '1321053'
$ws.Cells.Item(11,2).Value = 'https://aiesec.org/opportunity/global-talent/1321053'
$ws.Cells.Item(11,3).Value = 'International Sales Representetive German Speaker'
$ws.Cells.Item(11,4).Value = 'Maslak, Sarıyer/İstanbul, Türkiye'
$ws.Cells.Item(11,5).Value = 'No'
$ws.Cells.Item(11,6).Value = '17 applicants'
$ws.Cells.Item(11,7).Value = '6 - 18 Months'
$ws.Cells.Item(11,8).Value = 'Esvita Clinic'

# Row 12
$ws.Cells.Item(12,1).Value = '1321052'
$ws.Cells.Item(12,2).Value = 'https://aiesec.org/opportunity/global-talent/1321052'
$ws.Cells.Item(12,3).Value = 'International Sales Representetive'
$ws.Cells.Item(12,4).Value = 'Maslak, Sarıyer/İstanbul, Türkiye'
$ws.Cells.Item(12,5).Value = 'No'
$ws.Cells.Item(12,6).Value = '156 applicants'
$ws.Cells.Item(12,7).Value = '6 - 18 Months'
$ws.Cells.Item(12,8).Value = 'Esvita Clinic'

# Row 13
$ws.Cells.Item(13,1).Value = '1313793'
$ws.Cells.Item(13,2).Value = 'https://aiesec.org/opportunity/global-talent/1313793'
$ws.Cells.Item(13,3).Value = 'IT Sales Executive'
$ws.Cells.Item(13,4).Value = 'Kim Chung, Hoài Đức, Hà Nội, Việt Nam'
$ws.Cells.Item(13,5).Value = 'No'
$ws.Cells.Item(13,6).Value = '108 applicants'
$ws.Cells.Item(13,7).Value = '9 - 12 Weeks'
$ws.Cells.Item(13,8).Value = 'MOHA SOFTWARE JOINT STOCK COMPANY'

# Row 14
$ws.Cells.Item(14,1).Value = '1306716'
$ws.Cells.Item(14,2).Value = 'https://aiesec.org/opportunity/global-talent/1306716'
$ws.Cells.Item(14,3).Value = 'ENGINEERING'
$ws.Cells.Item(14,4).Value = 'Gaziantep, Türkiye'
$ws.Cells.Item(14,5).Value = 'No'
$ws.Cells.Item(14,6).Value = '28 applicants'
$ws.Cells.Item(14,7).Value = '6 - 18 Months'
$ws.Cells.Item(14,8).Value = 'Göymen Makarna'

# Row 15
$ws.Cells.Item(15,1).Value = '1289379'
$ws.Cells.Item(15,2).Value = 'https://aiesec.org/opportunity/global-talent/1289379'
$ws.Cells.Item(15,3).Value = 'Medical Advisor Portuguese Speaker'
$ws.Cells.Item(15,4).Value = 'İstanbul, Türkiye'
$ws.Cells.Item(15,5).Value = 'No'
$ws.Cells.Item(15,6).Value = '123 applicants'
$ws.Cells.Item(15,7).Value = '6 - 18 Months'
$ws.Cells.Item(15,8).Value = 'International Plus'

# Row 16
$ws.Cells.Item(16,1).Value = '1289378'
$ws.Cells.Item(16,2).Value = 'https://aiesec.org/opportunity/global-talent/1289378'
$ws.Cells.Item(16,3).Value = 'Medical Advisor (Spanish Speaker)'
$ws.Cells.Item(16,4).Value = 'İstanbul, Türkiye'
$ws.Cells.Item(16,5).Value = 'No'
$ws.Cells.Item(16,6).Value = '128 applicants'
$ws.Cells.Item(16,7).Value = '6 - 18 Months'
$ws.Cells.Item(16,8).Value = 'International Plus'

# Row 17
$ws.Cells.Item(17,1).Value = '1289377'
$ws.Cells.Item(17,2).Value = 'https://aiesec.org/opportunity/global-talent/1289377'
$ws.Cells.Item(17,3).Value = 'Medical Advisor (Italian Speaker)'
$ws.Cells.Item(17,4).Value = 'İstanbul, Türkiye'
$ws.Cells.Item(17,5).Value = 'No'
$ws.Cells.Item(17,6).Value = '40 applicants'
$ws.Cells.Item(17,7).Value = '6 - 18 Months'
$ws.Cells.Item(17,8).Value = 'International Plus'

# Row 18
$ws.Cells.Item(18,1).Value = '1289375'
$ws.Cells.Item(18,2).Value = 'https://aiesec.org/opportunity/global-talent/1289375'
$ws.Cells.Item(18,3).Value = 'Medical Advisor (German Speaker)'
$ws.Cells.Item(18,4).Value = 'İstanbul, Türkiye'
$ws.Cells.Item(18,5).Value = 'No'
$ws.Cells.Item(18,6).Value = '43 applicants'
$ws.Cells.Item(18,7).Value = '6 - 18 Months'
$ws.Cells.Item(18,8).Value = 'International Plus'
